$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP), G (sum)
# regenerated sval data to filter save games; F (Win) unchanged

$data = @{
    2  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 14.05633640148523)
    3  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    4  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    5  = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086)
    6  = @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086)
    7  = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    8  = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    9  = @(1.455362044514542, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 25.99591228164478)
    10 = @(1.455362044514542, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 25.99591228164478)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
